# Updated cryptos list on Thu Feb 29 04:54:03 UTC 2024 with GitHub Actions
#
# Applies per-row Price (D) / Volume(1h) (E) updates, and the two name/link
# swaps (rows 30<->31: Kaspa/LEO, rows 39<->40: LidoDAOToken/FirstDigitalUSD).
#
# Some Price values are plain decimals (e.g. "414.92") that Excel's COM layer
# would otherwise auto-convert to a Number on assignment. The source data is
# text (it mixes "."-grouped big numbers with plain decimals in the same
# column), so for those we force text by assigning with a leading apostrophe
# (the standard Excel "treat as text" marker) and then clear the resulting
# cell format so no stray number-format/quote-prefix style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force a value to be stored as text even if it looks like a number,
    # without leaving a lasting number-format/quote-prefix style on the cell.
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $text
        $range.ClearFormats()
    } else {
        $range.Value = $text
    }
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = '62.431.49'
$ws.Range("E2").Value = '  +9.25%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '3.447.22'
$ws.Range("E3").Value = '  +5.60%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  -0.01%  '

# Row 5 - BNB
Set-TextValue $ws.Range("D5") '414.92'
$ws.Range("E5").Value = '  +4.40%  '

# Row 6 - Solana
Set-TextValue $ws.Range("D6") '122.92'
$ws.Range("E6").Value = '  +13.05%  '

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = '3.443.63'
$ws.Range("E7").Value = '  +5.66%  '

# Row 8 - XRP
$ws.Range("E8").Value = '  +1.79%  '

# Row 9 - USDC
$ws.Range("E9").Value = '  +0.05%  '

# Row 10 - Cardano
$ws.Range("E10").Value = '  +4.65%  '

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") '0.128'
$ws.Range("E11").Value = '  +33.95%  '

# Row 12 - Avalanche
Set-TextValue $ws.Range("D12") '41.18'
$ws.Range("E12").Value = '  +4.21%  '

# Row 13 - TRON
$ws.Range("E13").Value = '  -0.31%  '

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '3.985.98'
$ws.Range("E14").Value = '  +5.45%  '

# Row 15 - Polkadot
Set-TextValue $ws.Range("D15") '8.48'
$ws.Range("E15").Value = '  +2.31%  '

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") '19.62'
$ws.Range("E16").Value = '  +3.16%  '

# Row 17 - WrappedEther
$ws.Range("D17").Value = '3.431.59'
$ws.Range("E17").Value = '  +4.82%  '

# Row 18 - WrappedBTC
$ws.Range("D18").Value = '62.223.29'
$ws.Range("E18").Value = '  +9.23%  '

# Row 19 - Polygon
$ws.Range("E19").Value = '  -0.64%  '

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") '10.78'
$ws.Range("E20").Value = '  -3.03%  '

# Row 21 - ShibaInu
$ws.Range("E21").Value = '  +22.24%  '

# Row 22 - ImmutableX
Set-TextValue $ws.Range("D22") '3.31'
$ws.Range("E22").Value = '  -1.19%  '

# Row 23 - Litecoin
Set-TextValue $ws.Range("D23") '81.77'
$ws.Range("E23").Value = '  +10.01%  '

# Row 24 - BitcoinCash
Set-TextValue $ws.Range("D24") '312.46'
$ws.Range("E24").Value = '  +6.67%  '

# Row 25 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D25") '12.95'
$ws.Range("E25").Value = '  -0.12%  '

# Row 26 - PancakeSwap
Set-TextValue $ws.Range("D26") '3.16'
$ws.Range("E26").Value = '  -0.66%  '

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") '31.05'
$ws.Range("E27").Value = '  +10.27%  '

# Row 28 - RenderToken
Set-TextValue $ws.Range("D28") '7.86'
$ws.Range("E28").Value = '  +5.72%  '

# Row 29 - Filecoin
Set-TextValue $ws.Range("D29") '7.75'
$ws.Range("E29").Value = '  -2.34%  '

# Row 30 - Kaspa->LEO (swap)
$ws.Range("B30").Value = 'LEO'
$ws.Range("C30").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range("D30") '4.30'
$ws.Range("E30").Value = '  -1.90%  '

# Row 31 - LEO->Kaspa (swap)
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D31") '0.174'
$ws.Range("E31").Value = '  +3.09%  '

# Row 32 - Hedera
$ws.Range("E32").Value = '  +3.73%  '

# Row 33 - Toncoin
Set-TextValue $ws.Range("D33") '2.57'
$ws.Range("E33").Value = '  +20.49%  '

# Row 34 - InjectiveProtocol
Set-TextValue $ws.Range("D34") '41.97'
$ws.Range("E34").Value = '  +4.62%  '

# Row 35 - Cosmos
Set-TextValue $ws.Range("D35") '11.36'
$ws.Range("E35").Value = '  +1.44%  '

# Row 36 - Dai
$ws.Range("E36").Value = '  +0.05%  '

# Row 37 - VeChain
$ws.Range("E37").Value = '  -1.64%  '

# Row 38 - OKB
Set-TextValue $ws.Range("D38") '52.69'
$ws.Range("E38").Value = '  +2.57%  '

# Row 39 - LidoDAOToken->FirstDigitalUSD (swap)
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D39") '0.996'
$ws.Range("E39").Value = '  -0.29%  '

# Row 40 - FirstDigitalUSD->LidoDAOToken (swap)
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D40") '3.50'
$ws.Range("E40").Value = '  +0.77%  '

# Row 41 - Stacks
Set-TextValue $ws.Range("D41") '3.06'
$ws.Range("E41").Value = '  +1.38%  '

# Row 42 - ARBITRUM
$ws.Range("E42").Value = '  +5.86%  '

# Row 43 - Stellar
$ws.Range("E43").Value = '  +2.75%  '

# Row 44 - Monero
Set-TextValue $ws.Range("D44") '134.81'
$ws.Range("E44").Value = '  -1.52%  '

# Row 45 - Celestia
Set-TextValue $ws.Range("D45") '17.08'
$ws.Range("E45").Value = '  +1.87%  '

# Row 46 - TheGraph
Set-TextValue $ws.Range("D46") '0.281'
$ws.Range("E46").Value = '  -0.98%  '

# Row 47 - NEARProtocol
Set-TextValue $ws.Range("D47") '3.87'
$ws.Range("E47").Value = '  -1.65%  '

# Row 48 - WEMIXToken
$ws.Range("E48").Value = '  -0.91%  '

# Row 49 - EnergySwap
Set-TextValue $ws.Range("D49") '21.91'
$ws.Range("E49").Value = '  -2.50%  '

# Row 50 - Maker
$ws.Range("D50").Value = '2.204.24'
$ws.Range("E50").Value = '  +2.35%  '

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = '3.780.24'
$ws.Range("E51").Value = '  +5.41%  '

